# RP3 ERT ATFM 2022 Jan-Dec — "Updates 2022 full year"
#
# 1. Bump the release date (ERT_ATFM_YY!B2); ERT_ATFM_MM / ERT_ATFM_FAB /
#    ERT_ATFM_LOC all read it back via `=ERT_ATFM_YY!B2`, so they refresh
#    automatically on recalculation.
# 2. On ERT_ATFM_LOC, the "Plan/Actual [2021]" column headers were stale
#    (the sheet actually reports 2022 numbers) — point them at the
#    existing "Plan [2022]" / "Actual [2022]" captions used elsewhere in
#    the workbook.
# 3. Refresh the per-entity breakdown: several ANSPs were renamed/merged
#    (Avinor -> Avinor Flysikring AS, DFS -> DFS + MUAC-DE, LPS -> LPS SR,
#    LVNL absorbs Maastricht UAC -> LVNL + MUAC-NL, NAV Portugal ->
#    NAV Portugal (Continental), skeyes splits into SE Oro Navigacija and
#    Belgium-Lux. + MUAC BE-LU, replacing the old Oro navigacija row), a
#    "Plan [2022]" percentage column B was populated, and the FLTS[TOT] /
#    delay-minute figures were refreshed. The table now has 28 entities
#    instead of 29, so the trailing row is blanked out.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Release date refresh
# ---------------------------------------------------------------------
$wsYY = $wb.Worksheets.Item("ERT_ATFM_YY")
$wsYY.Range("B2").Value = 45034

# ---------------------------------------------------------------------
# 2 & 3. ERT_ATFM_LOC entity table
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ERT_ATFM_LOC")

# Header row (row 5): correct the column captions to the 2022 figures
$ws.Range("B5").Value = "Plan [2022]"
$ws.Range("E5").Value = "Actual [2022]"

# Entity rows 6-33: Name (A), Plan [2022] share (B), FLTS [TOT] (C),
# En-route ATFM delay [min.] (D). Column E (share D/C) is a formula and
# recalculates on its own.
$entities = @(
    @{Row=6;  Name="ANS CR";                       Plan=0.11; C=550194;  D=824230},
    @{Row=7;  Name="Austro Control";                Plan=0.17; C=1049628; D=104902},
    @{Row=8;  Name="Avinor Flysikring AS";           Plan=0.08; C=529671;  D=3266},
    @{Row=9;  Name="BULATSA";                       Plan=0.08; C=822087;  D=50},
    @{Row=10; Name="Croatia Control";               Plan=0.16; C=712861;  D=407715},
    @{Row=11; Name="DCAC Cyprus";                   Plan=0.16; C=343957;  D=222},
    @{Row=12; Name="DFS + MUAC-DE";                 Plan=0.27; C=2516506; D=5728737},
    @{Row=13; Name="DSNA";                          Plan=0.25; C=2919041; D=3563319},
    @{Row=14; Name="EANS";                          Plan=0.03; C=141577;  D=80},
    @{Row=15; Name="ENAIRE";                        Plan=0.2;  C=1982636; D=675536},
    @{Row=16; Name="ENAV";                          Plan=0.11; C=1664146; D=362824},
    @{Row=17; Name="Fintraffic ANS";                Plan=0.05; C=192580;  D=0},
    @{Row=18; Name="HASP";                          Plan=0.14; C=895730;  D=138090},
    @{Row=19; Name="HungaroControl (EC)";           Plan=0.11; C=892471;  D=794061},
    @{Row=20; Name="IAA";                           Plan=0.03; C=582360;  D=603},
    @{Row=21; Name="LFV";                           Plan=0.07; C=541132;  D=22147},
    @{Row=22; Name="LGS";                           Plan=0.03; C=187979;  D=0},
    @{Row=23; Name="LPS SR";                        Plan=0.07; C=470158;  D=13907},
    @{Row=24; Name="LVNL + MUAC-NL";                Plan=0.14; C=1052844; D=49345},
    @{Row=25; Name="MATS";                          Plan=0.01; C=104143;  D=0},
    @{Row=26; Name="NAV Portugal (Continental)";    Plan=0.13; C=606631;  D=384482},
    @{Row=27; Name="NAVIAIR";                       Plan=0.06; C=516284;  D=762},
    @{Row=28; Name="PANSA";                         Plan=0.12; C=613073;  D=809805},
    @{Row=29; Name="ROMATSA";                       Plan=0.04; C=656325;  D=0},
    @{Row=30; Name="SE Oro Navigacija";             Plan=0.02; C=163350;  D=0},
    @{Row=31; Name="Belgium-Lux. + MUAC  BE-LU";    Plan=0.17; C=1037834; D=135538},
    @{Row=32; Name="Skyguide";                      Plan=0.19; C=1152707; D=392134},
    @{Row=33; Name="Slovenia Control";              Plan=0.09; C=360890;  D=156}
)

foreach ($e in $entities) {
    $r = $e.Row
    $ws.Cells.Item($r, 1).Value = $e.Name
    $ws.Cells.Item($r, 2).Value = $e.Plan
    $ws.Cells.Item($r, 3).Value = $e.C
    $ws.Cells.Item($r, 4).Value = $e.D
}

# Row 34 (formerly "Slovenia Control", now folded into row 33) is blanked.
$ws.Range("A34").ClearContents()
$ws.Range("B34").ClearContents()
$ws.Range("C34").ClearContents()
$ws.Range("D34").ClearContents()
$ws.Range("E34").ClearContents()
